# Updated cryptos list on Sat Oct 21 04:56:34 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on every touched cell so numeric-looking strings
# (prices like "0.0600" or "6.40") are preserved verbatim, matching the
# source data's inline-string representation instead of being reparsed as numbers.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.611.20"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +1.48%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.600.77"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +1.45%  "
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.51%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "212.43"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +0.53%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "26.82"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +3.61%  "
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  +1.36%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0600"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +1.19%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0911"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +0.99%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.828.23"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +1.32%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.602.39"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +1.18%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "29.618.44"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +1.54%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.537"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +3.02%  "
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +0.98%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.88"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "241.46"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  +1.23%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.62"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +2.30%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0694"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +0.42%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "3.98"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -0.18%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.23"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -0.68%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "155.08"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  +1.08%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "15.35"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +1.43%  "
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  +0.62%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "6.40"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +1.18%  "
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +0.46%  "
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  +2.77%  "
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +2.67%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.425.03"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell = $ws.Range("B35")
$cell.NumberFormat = "@"
$cell.Value = "MXToken"
$cell = $ws.Range("C35")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.89"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +5.31%  "
$cell = $ws.Range("B36")
$cell.NumberFormat = "@"
$cell.Value = "LidoDAOToken"
$cell = $ws.Range("C36")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.54"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +2.28%  "
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -2.21%  "
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  +0.31%  "
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +2.53%  "
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +2.88%  "
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +0.63%  "
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +4.95%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "54.19"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +1.23%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.808"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +2.51%  "
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +0.40%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.987"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +16.42%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "66.21"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +2.80%  "
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -0.80%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.740.65"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +1.41%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "85.97"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +0.22%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0104"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +2.99%  "
